$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with latest values.
# Numeric-looking price strings are forced back to Text so Excel keeps them as
# plain strings (matching the source data format) instead of auto-converting them
# to numbers.

$ws.Range("D2").Value = "41.933.85"
$ws.Range("E2").Value = "  +5.08%  "

$ws.Range("D3").Value = "2.254.36"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.94"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.55%  "

$ws.Range("E7").Value = "  +3.51%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.90%  "

$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("E14").Value = "  +3.30%  "

$ws.Range("D15").Value = "2.602.55"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").Value = "2.255.36"
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.757"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.33%  "

$ws.Range("D19").Value = "41.786.75"
$ws.Range("E19").Value = "  +4.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.84%  "

$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("E25").Value = "  +3.80%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.91%  "

$ws.Range("E29").Value = "  +12.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.08%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.94%  "

$ws.Range("E37").Value = "  +2.76%  "

$ws.Range("E38").Value = "  +3.52%  "

$ws.Range("E39").Value = "  +4.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.38%  "

$ws.Range("D43").Value = "2.050.51"
$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.93%  "

$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.68%  "

$ws.Range("E49").Value = "  +3.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.95%  "
